# Weekly fruit/vegetable price update: a new daily price record was added
# for "Achicoria" at Vega Modelo de Temuco, inserted as row 31. All the
# existing records from row 31 onward shift down by one row (to 32..110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31; rows 31-109 shift down to 32-110.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new price record.
$ws.Cells.Item(31, 1).Value2  = 10
$ws.Cells.Item(31, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value2  = "La Araucanía"
$ws.Cells.Item(31, 4).Value2  = 45076
$ws.Cells.Item(31, 5).Value2  = 9
$ws.Cells.Item(31, 6).Value2  = 100112010
$ws.Cells.Item(31, 7).Value2  = "Achicoria"
$ws.Cells.Item(31, 8).Value2  = "Sin especificar"
$ws.Cells.Item(31, 9).Value2  = "Primera"
$ws.Cells.Item(31, 10).Value2 = 85
$ws.Cells.Item(31, 11).Value2 = 10000
$ws.Cells.Item(31, 12).Value2 = 10000
$ws.Cells.Item(31, 13).Value2 = 10000
$ws.Cells.Item(31, 14).Value2 = "$/caja 18 unidades"
$ws.Cells.Item(31, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(31, 16).Value2 = 556
$ws.Cells.Item(31, 17).Value2 = 18
$ws.Cells.Item(31, 18).Value2 = "Hortaliza"
